# Week 5 - day 1
# Fix a units typo in the recitation problem: the density of air should be
# given in g/cm^3 (grams per cubic centimeter), not g/dm^3.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "1.29 g/dm^3 at 25",   # FindText
    $true,                 # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    1,                     # Wrap (wdFindContinue)
    $false,                # Format
    "1.29 g/cm^3 at 25",   # ReplaceWith
    2                      # Replace (wdReplaceAll)
)
